# Insert a new row at row 22 of the "general" sheet, adding the
# "errorPropTestEnable1" parameter row, shifting all rows below it down by 1.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("general")

# Insert a new blank row above the current row 22 (shifts 22..55 down to 23..56)
$ws.Rows("22:22").Insert()

# Populate the newly inserted row 22
$ws.Range("A22").Value = "errorPropTestEnable1"
$ws.Range("B22").Value = 1
$ws.Range("E22").Formula = "=B22"

# Match formatting/style of neighboring rows (same style ids as row 21 general rows)
$ws.Range("A22").Style = $ws.Range("A23").Style
$ws.Range("B22").Style = $ws.Range("B23").Style
$ws.Range("C22").Style = $ws.Range("C23").Style
$ws.Range("D22").Style = $ws.Range("D23").Style
$ws.Range("E22").Style = $ws.Range("E23").Style
